# Fruta / hortaliza, semanal
# Insert two new price records (rows 147-148) for "Early Burlat" cherries
# and two more new price records (rows 166-167, after the shift) for
# "Early Burlat" sold by 5-kilo trays, shifting the remaining weekly
# records down to keep the existing rows/data intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CerezaRow {
    param(
        [int]$Row,
        [double]$Fecha,
        [string]$Variedad,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [string]$Unidad,
        [string]$Origen,
        [double]$PrecioKg,
        [double]$KgUnidad
    )

    $ws.Cells.Item($Row, 1).Value = 3
    $ws.Cells.Item($Row, 2).Value = "Femacal de La Calera"
    $ws.Cells.Item($Row, 3).Value = "Coquimbo"
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = 5
    $ws.Cells.Item($Row, 6).Value = "Fruta"
    $ws.Cells.Item($Row, 7).Value = 100103
    $ws.Cells.Item($Row, 8).Value = "Frutos de hueso (carozo)"
    $ws.Cells.Item($Row, 9).Value = 100103001
    $ws.Cells.Item($Row, 10).Value = "Cereza"
    $ws.Cells.Item($Row, 11).Value = $Variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $Unidad
    $ws.Cells.Item($Row, 18).Value = $Origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $KgUnidad
}

# --- First insertion point: two new rows before the old row 147 ---
$ws.Rows.Item(147).Insert()
$ws.Rows.Item(147).Insert()

Set-CerezaRow 147 44508 "Early Burlat" "Primera" 48 50000 50000 50000 "$/bandeja 10 kilos" "Provincia de Curicó" 5000 10
Set-CerezaRow 148 44508 "Early Burlat" "Segunda" 50 40000 40000 40000 "$/bandeja 10 kilos" "Provincia de Curicó" 4000 10

# --- Second insertion point: two new rows before what is now row 166 ---
$ws.Rows.Item(166).Insert()
$ws.Rows.Item(166).Insert()

Set-CerezaRow 166 44504 "Early Burlat" "Primera" 56 27000 27000 27000 "$/bandeja 5 kilos" "Provincia de Curicó" 5400 5
Set-CerezaRow 167 44505 "Early Burlat" "Primera" 38 27000 27000 27000 "$/bandeja 5 kilos" "Provincia de Curicó" 5400 5
